$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text changes (row1 = CN description, row2 = type, row3 = field name) ---
# Column E: HeroSkillLevel -> CardLevel ("int" -> "double"; "英雄技能等级" -> "怪物卡等级")
$ws.Range("E1").Value = "怪物卡等级"
$ws.Range("E2").Value = "double"
$ws.Range("E3").Value = "CardLevel"

# --- Column B (Exp) tweaks: B6 literal bump, then B7+ becomes a formula chain ---
$ws.Range("B6").Value = 50
$ws.Range("B7").Formula = "=INT(B6*0.7+A7*10+A7*A7*3-30)"
$ws.Range("B8:B71").Formula = "=INT(B7*0.7+A8*10+A8*A8*3-30)"
$ws.Range("B72:B102").Formula = "=INT(B71*0.7+A72*10+A72*A72*3-30)"

# --- Column D (TowerLevel) becomes a calculated column: INT(E) ---
$ws.Range("D4").Formula = "=INT(E4)"
$ws.Range("D5:D68").Formula = "=INT(E5)"
$ws.Range("D69:D102").Formula = "=INT(E69)"

# --- Column E (CardLevel) values: literal ramp 1.0 -> 2.9 for rows 4-23, then formula-driven ---
$ws.Range("E5").Value = 1.1
$ws.Range("E6").Value = 1.2
$ws.Range("E7").Value = 1.3
$ws.Range("E8").Value = 1.4
$ws.Range("E9").Value = 1.5
$ws.Range("E10").Value = 1.6
$ws.Range("E11").Value = 1.7
$ws.Range("E12").Value = 1.8
$ws.Range("E13").Value = 1.9
$ws.Range("E15").Value = 2.1
$ws.Range("E16").Value = 2.2
$ws.Range("E17").Value = 2.3
$ws.Range("E18").Value = 2.4
$ws.Range("E19").Value = 2.5
$ws.Range("E20").Value = 2.6
$ws.Range("E21").Value = 2.7
$ws.Range("E22").Value = 2.8
$ws.Range("E23").Value = 2.9

$ws.Range("E24").Formula = "=3.3+(A24-21)*0.1-E23*0.12"
$ws.Range("E25:E88").Formula = "=3.3+(A25-21)*0.1-E24*0.12"
$ws.Range("E89:E102").Formula = "=3.3+(A89-21)*0.1-E88*0.12"

# --- Column widths: split D/E into separate <col> entries, E grows to width 8 ---
$ws.Columns.Item(5).ColumnWidth = 7.2857142857142865

# --- Selection shown when the sheet is active ---
$ws.Range("B7:B102").Select()
